$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 184.0626906666667
$ws.Range("H2").Value = 552.188072
$ws.Range("I2").Value = 0.6510505751503485
$ws.Range("J2").Value = 0.6510505751503486
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8383893333333333
$ws.Range("N2").Value = 2.515168
$ws.Range("Q2").Value = 154.3161965195662
$ws.Range("R2").Value = 1388.845768676096
$ws.Range("S2").Value = 0.6510505751503485
$ws.Range("T2").Value = 0.6510505751503486

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 57.4434
$ws.Range("H3").Value = 172.3302
$ws.Range("I3").Value = 0.2031838091312023
$ws.Range("J3").Value = 0.2031838091312023
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8383893333333333
$ws.Range("N3").Value = 2.515168
$ws.Range("Q3").Value = 48.1599338304
$ws.Range("R3").Value = 433.4394044736
$ws.Range("S3").Value = 0.2031838091312023
$ws.Range("T3").Value = 0.2031838091312023

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 41.21033366666666
$ws.Range("H4").Value = 123.631001
$ws.Range("I4").Value = 0.1457656157184491
$ws.Range("J4").Value = 0.1457656157184491
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8383893333333333
$ws.Range("N4").Value = 2.515168
$ws.Range("Q4").Value = 34.55030416924089
$ws.Range("R4").Value = 310.952737523168
$ws.Range("S4").Value = 0.1457656157184491
$ws.Range("T4").Value = 0.1457656157184491
